# Move humidity set points from code to construction property database
# (INDOOR_COMFORT sheet): add rhum_min_pc / rhum_max_pc columns (G, H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INDOOR_COMFORT")

# --- New header cells (G1, H1) ---
$ws.Range("G1").Value = "rhum_min_pc"
$ws.Range("H1").Value = "rhum_max_pc"

# --- New data columns: every data row (2-20) gets 30 / 70 ---
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 7).Value = 30
    $ws.Cells.Item($r, 8).Value = 70
}

# --- Match formatting of existing columns ---
# Header style (same as the other header cells, e.g. A1)
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Data style (same as the already-correctly-styled F2 cell)
$ws.Range("F2").Copy()
$ws.Range("G2:H20").PasteSpecial(-4122)

# Column F (Ve_lps) used an inconsistent style (s=8) on most rows;
# normalize it to match the rest (s=5), same as F2/F3 already have.
$ws.Range("F4:F20").PasteSpecial(-4122)

# --- Active sheet / selection bookkeeping ---
# INDOOR_COMFORT becomes the active/selected tab (was SUPPLY).
$ws.Activate()
$ws.Range("H2").Select()
